$wb = $excel.ActiveWorkbook

# 1. Rename the original sheet "Original" -> "for BC"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "for BC"

# 2. Add a new "cheat sheet" worksheet right after "for BC"
$cheat = $wb.Worksheets.Add($null, $ws)
$cheat.Name = "cheat sheet"
$cheat.Range("A1").Value = "Detected"
$cheat.Range("A2").Value = "Not Detected"
$cheat.Range("A3").Value = "N/A"
$cheat.Columns.Item(1).ColumnWidth = 12.85546875

# 3. Point the existing list data validation (on B12:C54 of "for BC") at the
#    new cheat sheet list instead of the broken #REF! reference.
$rng = $ws.Range("B12:C54")
$rng.Validation.Delete() | Out-Null
$rng.Validation.Add(3, 1, 1, "='cheat sheet'!`$A`$1:`$A`$3") | Out-Null
$rng.Validation.IgnoreBlank = $true
$rng.Validation.InCellDropdown = $true
$rng.Validation.ShowInput = $true
$rng.Validation.ShowError = $true

# 4. Collapse the sheet selection from B12:B54 down to just B12
$ws.Range("B12").Select() | Out-Null
$ws.Activate() | Out-Null
